$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - values shifted: B2 removed, C2/D2/E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 0.54866657396459306
$ws.Range("D2").Value = 0.18420912400535153
$ws.Range("E2").Value = 1.0593207285909152

# Row 3 - values shifted: C3 removed, D3 added, B3/E3 updated
$ws.Range("B3").Value = 0.40603904931766521
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 1.0845369198341224
$ws.Range("E3").Value = 1.6539524082197778

# Update selection to match new active range
$ws.Range("B1:E3").Select()
